$wb = $excel.ActiveWorkbook

# The workbook originally has a single sheet named "Sheet1" containing the
# species/synonyms table. Worksheets.Add() with no args inserts a new,
# blank sheet immediately before the currently active sheet (matching
# Excel's default "Insert Sheet" behaviour) and selects it; since the name
# "Sheet1" is already taken, the new sheet is auto-named "Sheet2" while the
# original sheet keeps its name "Sheet1" and both its data and its
# (now de-selected) tab state.
$newSheet = $wb.Worksheets.Add()

$newSheet.Range("A1").Value = "New Phytologist Supporting Information"
$newSheet.Range("A2").Value = "Photographs as an essential biodiversity resource: drivers of gaps in the vascular plant photographic record"
$newSheet.Range("A3").Value = "Thomas Mesaglio, Hervé Sauquet, David Coleman, Elizabeth Wenk, William K Cornwell"
$newSheet.Range("A4").Value = "Accepted 8 February 2023"
$newSheet.Range("A6").Value = "Caption"
$newSheet.Range("A7").Value = "List of all species for which photographs were located under a name different to that currently accepted by the Australian Plant Census (synonym, orthographic variant, etc.). "

$newSheet.Range("A1").Font.Bold = $true
$newSheet.Range("A6").Font.Bold = $true

[void]$newSheet.Range("I15").Select()
